$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 90
$ws.Range("F5").Value = 1207
$ws.Range("F6").Value = 9446
$ws.Range("F7").Value = 7275
$ws.Range("F8").Value = 198
$ws.Range("F9").Value = 327
$ws.Range("F10").Value = 6063
$ws.Range("F12").Value = 81
$ws.Range("F13").Value = 30
$ws.Range("F14").Value = 6713
$ws.Range("F15").Value = 1120
$ws.Range("F16").Value = 481
$ws.Range("F17").Value = 455
$ws.Range("F18").Value = 34
$ws.Range("F19").Value = 659
$ws.Range("F21").Value = 298
$ws.Range("F25").Value = 10899
$ws.Range("F27").Value = 51
$ws.Range("F28").Value = 2070
$ws.Range("F29").Value = 2659
$ws.Range("F30").Value = 48
$ws.Range("F32").Value = 2415
$ws.Range("F33").Value = 93
$ws.Range("F35").Value = 32
$ws.Range("F38").Value = 1503
$ws.Range("F40").Value = 29
$ws.Range("F41").Value = 5529
$ws.Range("F42").Value = 1227
$ws.Range("F43").Value = 783
$ws.Range("F44").Value = 142
$ws.Range("F46").Value = 1090
$ws.Range("F47").Value = 1446
$ws.Range("F48").Value = 80
$ws.Range("F49").Value = 1113

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 3
$ws.Range("F8").Value = 245
$ws.Range("F11").Value = 204
$ws.Range("F12").Value = 2
$ws.Range("F20").Value = 39
$ws.Range("F21").Value = 9

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 119
$ws.Range("F3").Value = 219

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 90
$ws.Range("F4").Value = 1207
$ws.Range("F5").Value = 9446
$ws.Range("F6").Value = 7275
$ws.Range("F7").Value = 119
$ws.Range("F8").Value = 198
$ws.Range("F9").Value = 219
$ws.Range("F12").Value = 6064
$ws.Range("F13").Value = 81
$ws.Range("F15").Value = 30
$ws.Range("F16").Value = 6713
$ws.Range("F17").Value = 6713
$ws.Range("F18").Value = 1120
$ws.Range("F19").Value = 481
$ws.Range("F20").Value = 455
$ws.Range("F21").Value = 659
$ws.Range("F23").Value = 298
$ws.Range("F25").Value = 245
$ws.Range("F27").Value = 204
$ws.Range("F28").Value = 10899
$ws.Range("F30").Value = 51
$ws.Range("F31").Value = 2070
$ws.Range("F32").Value = 2659
$ws.Range("F33").Value = 2415
$ws.Range("F34").Value = 93
$ws.Range("F36").Value = 32
$ws.Range("F39").Value = 1503
$ws.Range("F40").Value = 5529
$ws.Range("F41").Value = 39
$ws.Range("F42").Value = 1227
$ws.Range("F43").Value = 783
$ws.Range("F44").Value = 142
$ws.Range("F46").Value = 1090
$ws.Range("F48").Value = 1446
$ws.Range("F49").Value = 80
$ws.Range("F50").Value = 1113
